# Add 2022-Q3 data.
#
# Before: Worksheets = [ "总计", "2021-Q3" ]
# After : Worksheets = [ "总计", "2022-Q3", "2021-Q3" ]
#   - a new "2022-Q3" sheet (holding fund-holding detail rows) is inserted
#     right before the existing "2021-Q3" sheet
#   - the "总计" (totals) sheet gets a new row for 2022-Q3 (inserted above
#     the existing 2021-Q3 totals row)

$wb = $excel.ActiveWorkbook

$wsTotal   = $wb.Worksheets.Item(1)   # "总计"
$wsOld2021 = $wb.Worksheets.Item(2)   # "2021-Q3" (existing fund detail sheet)

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet immediately before "2021-Q3" so
#    the final tab order is 总计, 2022-Q3, 2021-Q3.
# ---------------------------------------------------------------------
$wsNew = $wb.Worksheets.Add($wsOld2021)
$wsNew.Name = "2022-Q3"

# Match the page setup used by the rest of the workbook.
$wsNew.PageSetup.LeftMargin   = 54
$wsNew.PageSetup.RightMargin  = 54
$wsNew.PageSetup.TopMargin    = 72
$wsNew.PageSetup.BottomMargin = 72
$wsNew.PageSetup.HeaderMargin = 36
$wsNew.PageSetup.FooterMargin = 36
$wsNew.Outline.SummaryRow    = 1
$wsNew.Outline.SummaryColumn = 1

# ---------------------------------------------------------------------
# 2. "总计" sheet: the 2022-Q3 totals take over row 2, and the old
#    2021-Q3 totals (previously row 2) move down to row 3.
# ---------------------------------------------------------------------
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Range("A3").Value2 = 1
$wsTotal.Range("B3").Value2 = "2021-Q3"
$wsTotal.Range("C3").Value2 = 2
$wsTotal.Range("D3").Value2 = 0.13

$wsTotal.Range("A2").Value2 = 0
$wsTotal.Range("B2").Value2 = "2022-Q3"
$wsTotal.Range("C2").Value2 = 1
$wsTotal.Range("D2").Value2 = 0.18

# ---------------------------------------------------------------------
# 3. Populate the new "2022-Q3" sheet: header row (same look as the
#    other sheets' header rows) plus one fund holding data row.
# ---------------------------------------------------------------------
$wsTotal.Range("B1").Copy($wsNew.Range("B1"))
$wsTotal.Range("B1").Copy($wsNew.Range("C1"))
$wsTotal.Range("B1").Copy($wsNew.Range("D1"))
$wsTotal.Range("B1").Copy($wsNew.Range("E1"))
$wsTotal.Range("B1").Copy($wsNew.Range("F1"))
$wsTotal.Range("B1").Copy($wsNew.Range("G1"))
$wsTotal.Range("B1").Copy($wsNew.Range("H1"))

$wsNew.Range("B1").Value2 = "基金代码"
$wsNew.Range("C1").Value2 = "基金名称"
$wsNew.Range("D1").Value2 = "基金规模"
$wsNew.Range("E1").Value2 = "股票总仓位"
$wsNew.Range("F1").Value2 = "仓位占比"
$wsNew.Range("G1").Value2 = "持有市值(亿元)"
$wsNew.Range("H1").Value2 = "仓位排名"

$wsTotal.Range("A2").Copy($wsNew.Range("A2"))
$wsNew.Range("A2").Value2 = 0
# Fund code / size / positions are stored as text (matching the sibling
# "2021-Q3" sheet), so force text with a leading apostrophe for the
# numeric-looking values.
$wsNew.Range("B2").Value2 = "'513360"
$wsNew.Range("C2").Value2 = "博时中证全球中国教育主题ETF（QDII）"
$wsNew.Range("D2").Value2 = "'4.81"
$wsNew.Range("E2").Value2 = "'99.43"
$wsNew.Range("F2").Value2 = "'3.78"
$wsNew.Range("G2").Value2 = "'0.1818"
$wsNew.Range("H2").Value2 = 6
